# Saldo_guide.xlsx update
# - Advance every "Dt. Referencia" (column G, rows 2:310) from 2024-04-03 (45385)
#   to 2024-04-04 (45386).
# - Refresh the "Saldo Previsto" (D) / "Vl. Total" (H) figures for the accounts
#   whose balances moved since the last extract.
# - Leave the active selection on O26, matching where the author's cursor
#   ended up when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump every reference date in column G (rows 2-310) by one day.
$ws.Range("G2:G310").Value = 45386

# Updated balances (Saldo Previsto / Vl. Total columns) for the rows whose
# totals changed in this extract.
$updates = @{
    23  = 23008.2
    47  = 18983.21
    55  = 34014.33
    63  = 5980.53
    73  = 2009.29
    115 = 39991.27
    118 = 30917.07
    120 = 96957.97
    121 = 11030.06
    125 = 29865.85
    126 = 27850.74
    129 = 68000
    130 = 0
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $value
    $ws.Cells.Item($row, 8).Value = $value
}

# Restore the cursor position recorded by the author at save time.
$ws.Range("O26").Select() | Out-Null
